$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Part 1: For a handful of match rows, the data that was entered for the
# "home" row and the "away" row got swapped during the source re-scrape.
# Everything except column A (the sequential id, which stays tied to the
# row) needs to be exchanged between each pair of rows. Column E (the match
# Date) is identical for both rows of a pair, so it is skipped to avoid
# needlessly rewriting the date cell (and is excluded from the swap ranges).
# ---------------------------------------------------------------------------

function Swap-RowData($r1, $r2) {
    $rangeA1  = "B$r1" + ":D$r1"
    $rangeB1  = "B$r2" + ":D$r2"
    $rangeA2  = "F$r1" + ":AC$r1"
    $rangeB2  = "F$r2" + ":AC$r2"
    $scratch1 = "B1000:D1000"
    $scratch2 = "F1000:AC1000"

    $ws.Range($rangeA1).Copy()
    $ws.Range($scratch1).PasteSpecial(-4104)
    $ws.Range($rangeB1).Copy()
    $ws.Range($rangeA1).PasteSpecial(-4104)
    $ws.Range($scratch1).Copy()
    $ws.Range($rangeB1).PasteSpecial(-4104)
    $ws.Range($scratch1).Clear()

    $ws.Range($rangeA2).Copy()
    $ws.Range($scratch2).PasteSpecial(-4104)
    $ws.Range($rangeB2).Copy()
    $ws.Range($rangeA2).PasteSpecial(-4104)
    $ws.Range($scratch2).Copy()
    $ws.Range($rangeB2).PasteSpecial(-4104)
    $ws.Range($scratch2).Clear()
}

Swap-RowData 47 48
Swap-RowData 109 110
Swap-RowData 149 150
Swap-RowData 229 230
Swap-RowData 232 233
Swap-RowData 245 246
Swap-RowData 263 264

# ---------------------------------------------------------------------------
# Part 2: Two brand-new matches were appended at the bottom of the sheet.
# Copy the formatting from the last existing row and then fill in the data.
# These rows have no FTHG/FTAG/FTR/PL_Ahh/PL_AhUnder yet (match not played).
# ---------------------------------------------------------------------------

$ws.Range("A295:AC295").Copy()
$ws.Range("A296:AC296").PasteSpecial(-4122)
$ws.Cells.Item(296,1).Value = 294
$ws.Cells.Item(296,2).Value = 7645829
$ws.Cells.Item(296,3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(296,4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(296,5).Value = 45403.625
$ws.Cells.Item(296,6).Value = "Puebla Women"
$ws.Cells.Item(296,7).Value = "Juarez FC Women"
$ws.Cells.Item(296,11).Value = 3.4
$ws.Cells.Item(296,12).Value = 3.75
$ws.Cells.Item(296,13).Value = 1.8
$ws.Cells.Item(296,14).Value = 4.75
$ws.Cells.Item(296,15).Value = 4.2
$ws.Cells.Item(296,16).Value = 1.5
$ws.Cells.Item(296,17).Value = 1
$ws.Cells.Item(296,18).Value = 1.9
$ws.Cells.Item(296,19).Value = 1.9
$ws.Cells.Item(296,20).Value = 3
$ws.Cells.Item(296,21).Value = 1.925
$ws.Cells.Item(296,22).Value = 1.875
$ws.Cells.Item(296,23).Value = 0
$ws.Cells.Item(296,24).Value = 0
$ws.Cells.Item(296,25).Value = 0
$ws.Cells.Item(296,26).Value = 0
$ws.Cells.Item(296,27).Value = 0
$ws.Range("H296:J296").ClearContents()
$ws.Range("AB296:AC296").ClearContents()

$ws.Range("A295:AC295").Copy()
$ws.Range("A297:AC297").PasteSpecial(-4122)
$ws.Cells.Item(297,1).Value = 295
$ws.Cells.Item(297,2).Value = 7645830
$ws.Cells.Item(297,3).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(297,4).Value = "Mexico Liga MX Femenil"
$ws.Cells.Item(297,5).Value = 45403.79166666666
$ws.Cells.Item(297,6).Value = "Pachuca Women"
$ws.Cells.Item(297,7).Value = "Chivas Guadalajara Women"
$ws.Cells.Item(297,11).Value = 1.727
$ws.Cells.Item(297,12).Value = 3.6
$ws.Cells.Item(297,13).Value = 4
$ws.Cells.Item(297,14).Value = 2.1
$ws.Cells.Item(297,15).Value = 3.5
$ws.Cells.Item(297,16).Value = 2.8
$ws.Cells.Item(297,17).Value = -0.25
$ws.Cells.Item(297,18).Value = 1.95
$ws.Cells.Item(297,19).Value = 1.85
$ws.Cells.Item(297,20).Value = 3
$ws.Cells.Item(297,21).Value = 1.875
$ws.Cells.Item(297,22).Value = 1.925
$ws.Cells.Item(297,23).Value = 0
$ws.Cells.Item(297,24).Value = 0
$ws.Cells.Item(297,25).Value = 0
$ws.Cells.Item(297,26).Value = 0
$ws.Cells.Item(297,27).Value = 0
$ws.Range("H297:J297").ClearContents()
$ws.Range("AB297:AC297").ClearContents()
